$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.032.80'
$ws.Range("E2").Value = '  -4.80%  '
$ws.Range("D3").Value = '2.921.05'
$ws.Range("E3").Value = '  -7.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '476.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -9.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.21%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '2.915.95'
$ws.Range("E8").Value = '  -7.78%  '
$ws.Range("E9").Value = '  -8.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.77'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0985'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -11.49%  '
$ws.Range("E12").Value = '  -13.09%  '
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("D14").Value = '3.426.74'
$ws.Range("E14").Value = '  -7.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.33%  '
$ws.Range("D16").Value = '55.010.43'
$ws.Range("E16").Value = '  -4.81%  '
$ws.Range("D17").Value = '2.919.04'
$ws.Range("E17").Value = '  -7.03%  '
$ws.Range("E18").Value = '  -11.48%  '
$ws.Range("E19").Value = '  -5.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -11.35%  '
$ws.Range("E21").Value = '  -10.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '305.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -12.32%  '
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.450'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -12.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '59.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -14.94%  '
$ws.Range("E27").Value = '  -7.37%  '
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").Value = '0.0₃0821'
$ws.Range("E29").Value = '  -15.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.45%  '
$ws.Range("E33").Value = '  -12.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.99'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -12.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '145.69'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -14.02%  '
$ws.Range("E37").Value = '  -12.57%  '
$ws.Range("E38").Value = '  -12.51%  '
$ws.Range("E39").Value = '  -10.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0628'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.99%  '
$ws.Range("D41").Value = '2.951.23'
$ws.Range("E41").Value = '  -7.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.78'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -11.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.974'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -10.43%  '
$ws.Range("E45").Value = '  -10.90%  '
$ws.Range("E46").Value = '  -8.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -12.53%  '
$ws.Range("D48").Value = '2.072.94'
$ws.Range("E48").Value = '  -8.67%  '
$ws.Range("E49").Value = '  -12.47%  '
$ws.Range("E50").Value = '  -6.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -11.50%  '
